$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.074.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.71"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5240"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2596"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06293"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.66"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07662"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.626.62"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.411"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.861.51"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5519"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8272"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.91"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.051.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.684"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.27"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.155"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.58"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1216"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.51%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.395"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05965"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.254"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.437"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.404"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.640"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9830"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.394"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.760"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5671"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01616"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.733"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.032.59"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.11"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.787.16"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.74"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9959"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.046"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05153"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4215"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.68%  "
